$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells for rows 2-42 (ticker symbols shifted/replaced)
$ws.Range("B2").Value = "NSE:AETHER"
$ws.Range("C2").Value = "NSE:ALANKIT"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "NSE:BERGEPAINT"
$ws.Range("F2").Value = "NSE:GRASIM"
$ws.Range("B3").Value = "NSE:ARVINDFASN"
$ws.Range("C3").Value = "NSE:ALPHAGEO"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "NSE:CIPLA"
$ws.Range("F3").Value = "NSE:ICICIBANK"
$ws.Range("B4").Value = "NSE:ASHIMASYN"
$ws.Range("C4").Value = "NSE:ARTNIRMAN"
$ws.Range("E4").Value = "NSE:COALINDIA"
$ws.Range("F4").Value = "NSE:MOTHERSON"
$ws.Range("B5").Value = "NSE:CHOICEIN"
$ws.Range("C5").Value = "NSE:ASAHISONG"
$ws.Range("B6").Value = "NSE:COCHINSHIP"
$ws.Range("C6").Value = "NSE:ASTERDM"
$ws.Range("B7").Value = "NSE:EMSLIMITED"
$ws.Range("C7").Value = "NSE:CARTRADE"
$ws.Range("B8").Value = "NSE:GICHSGFIN"
$ws.Range("C8").Value = "NSE:CENTENKA"
$ws.Range("B9").Value = "NSE:GRASIM"
$ws.Range("C9").Value = "NSE:CONFIPET"
$ws.Range("B10").Value = "NSE:GULFOILLUB"
$ws.Range("C10").Value = "NSE:CREDITACC"
$ws.Range("B11").Value = "NSE:HEXATRADEX"
$ws.Range("C11").Value = "NSE:EIHAHOTELS"
$ws.Range("B12").Value = "NSE:KAKATCEM"
$ws.Range("C12").Value = "NSE:GENESYS"
$ws.Range("B13").Value = "NSE:KAMOPAINTS"
$ws.Range("C13").Value = "NSE:GICRE"
$ws.Range("B14").Value = "NSE:KESORAMIND"
$ws.Range("C14").Value = "NSE:GRAVITA"
$ws.Range("B15").Value = "NSE:MAZDOCK"
$ws.Range("C15").Value = "NSE:HERCULES"
$ws.Range("B16").Value = "NSE:MBLINFRA"
$ws.Range("C16").Value = "NSE:HEUBACHIND"
$ws.Range("B17").Value = "NSE:MTNL"
$ws.Range("C17").Value = "NSE:HGS"
$ws.Range("B18").Value = "NSE:NAGREEKEXP"
$ws.Range("C18").Value = "NSE:HILTON"
$ws.Range("B19").Value = "NSE:NRBBEARING"
$ws.Range("C19").Value = "NSE:HINDPETRO"
$ws.Range("B20").Value = "NSE:PAVNAIND"
$ws.Range("C20").Value = "NSE:INDIANHUME"
$ws.Range("B21").Value = "NSE:PTC"
$ws.Range("C21").Value = "NSE:J&KBANK"
$ws.Range("B22").Value = "NSE:RAYMOND"
$ws.Range("C22").Value = "NSE:JINDALSTEL"
$ws.Range("B23").Value = "NSE:RSWM"
$ws.Range("C23").Value = "NSE:KANORICHEM"
$ws.Range("B24").Value = "NSE:SALZERELEC"
$ws.Range("C24").Value = "NSE:KOTHARIPET"
$ws.Range("C25").Value = "NSE:KTKBANK"
$ws.Range("C26").Value = "NSE:LINC"
$ws.Range("C27").Value = "NSE:MAITHANALL"
$ws.Range("C28").Value = "NSE:MARATHON"
$ws.Range("C29").Value = "NSE:MASTEK"
$ws.Range("C30").Value = "NSE:MEGASOFT"
$ws.Range("C31").Value = "NSE:MOL"
$ws.Range("C32").Value = "NSE:MTARTECH"
$ws.Range("C33").Value = "NSE:NACLIND"
$ws.Range("C34").Value = "NSE:NATCOPHARM"
$ws.Range("C35").Value = "NSE:NBCC"
$ws.Range("C36").Value = "NSE:NESCO"
$ws.Range("C37").Value = "NSE:ORIENTHOT"
$ws.Range("C38").Value = "NSE:PRICOLLTD"
$ws.Range("C39").Value = "NSE:PRITIKAUTO"
$ws.Range("C40").Value = "NSE:RALLIS"
$ws.Range("C41").Value = "NSE:RAMKY"
$ws.Range("C42").Value = "NSE:RITES"

# Remove now-unused trailing rows 43-50 (table shrank from 48 to 41 data rows)
$ws.Range("A43:F50").EntireRow.Delete()
